# Apply "Added all mandetory debug functionallity + tests of FPGA" edit
# to the DebugDataTransfer sheet.
#
# Summary of the change:
#  - Row 52: add a new "BankID 4 Bit" cell in C52, and rename the
#    "read_iram" description in E52 to "read_ram".
#  - Insert two brand-new documentation rows (53 and 54) describing the
#    new 0x33 "VRAM Write" and 0x34 "MMIO Write" memory instructions.
#    (These reuse previously-empty row numbers, so nothing below shifts.)
#  - Move the active selection to D17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52: existing "Memory" block row for 0x32 -----------------------
# New cell: BankID 4 Bit
$ws.Range("C52").Value = "BankID 4 Bit"
# Renamed description: read_iram -> read_ram
$ws.Range("E52").Value = "read_ram"

# --- Row 53 (new): 0x33 VRAM Write ---------------------------------------
$ws.Range("B53").Value = "0x33"
$ws.Range("C53").Value = "Addr"
$ws.Range("D53").Value = "Data"
$ws.Range("E53").Value = "VRAM Write"
$ws.Range("G53").Value = "0x33"

# --- Row 54 (new): 0x34 MMIO Write ---------------------------------------
$ws.Range("B54").Value = "0x34"
$ws.Range("C54").Value = "Addr"
$ws.Range("D54").Value = "Data"
$ws.Range("E54").Value = "MMIO Write"
$ws.Range("G54").Value = "0x34"

# --- Update the selected/active cell -------------------------------------
$ws.Range("D17").Select()
